$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "last row date" number format (style 3) before we move it
$lastRowDateFormat = $ws.Range("S11").NumberFormat

# Populate new row 12 with the latest bunker price data
$ws.Range("A12").Value = 528
$ws.Range("B12").Value = 570
$ws.Range("C12").Value = 542
$ws.Range("D12").Value = 518
$ws.Range("E12").Value = 651
$ws.Range("F12").Value = 525
$ws.Range("G12").Value = 630
$ws.Range("H12").Value = 661
$ws.Range("I12").Value = 572
$ws.Range("J12").Value = 526
$ws.Range("K12").Value = 567
$ws.Range("L12").Value = 513
$ws.Range("M12").Value = 583
$ws.Range("N12").Value = 521
$ws.Range("O12").Value = 650
$ws.Range("P12").Value = 767
$ws.Range("Q12").Value = 560
$ws.Range("R12").Value = 667
$ws.Range("T12").Value = 590
$ws.Range("U12").Value = 596
$ws.Range("V12").Value = 608
$ws.Range("W12").Value = 508
$ws.Range("X12").Value = 513
$ws.Range("Y12").Value = 545
$ws.Range("Z12").Value = 757
$ws.Range("AA12").Value = 537
$ws.Range("AB12").Value = 575
$ws.Range("AC12").Value = 526
$ws.Range("AD12").Value = 652
$ws.Range("AE12").Value = 604.5
$ws.Range("AF12").Value = 563
$ws.Range("AG12").Value = 530
$ws.Range("AH12").Value = 575
$ws.Range("AI12").Value = 883
$ws.Range("AJ12").Value = 650
$ws.Range("AK12").Value = 510
$ws.Range("AL12").Value = 625
$ws.Range("AM12").Value = 551
$ws.Range("AN12").Value = 517
$ws.Range("AO12").Value = 535
$ws.Range("AP12").Value = 513
$ws.Range("AQ12").Value = 508
$ws.Range("AR12").Value = 497
$ws.Range("AS12").Value = 525
$ws.Range("AT12").Value = 550
$ws.Range("AU12").Value = 497
$ws.Range("AV12").Value = 560

# New date for row 12
$ws.Range("S12").Value = 45742

# S11 switches to the regular datetime format used by earlier rows (style 2)
$ws.Range("S11").NumberFormat = $ws.Range("S2").NumberFormat

# S12 (new last row) takes on the "date only" format previously on S11 (style 3)
$ws.Range("S12").NumberFormat = $lastRowDateFormat
